$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Columns.Item(3).Insert()
$ws.Range("C1").Value = "funding_source"
$ws.Range("C4").Value = "NC"
$ws.Range("A4:C4").Font.Size = 11
